$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.239.12'
$ws.Range('E2').Value = '  +0.46%  '

# Row 3
$ws.Range('D3').Value = '1.844.70'
$ws.Range('E3').Value = '  +0.40%  '

# Row 4
$ws.Range('E4').Value = '  +0.15%  '

# Row 5
$ws.Range('D5').Value = "'244.39"
$ws.Range('E5').Value = '  -0.18%  '

# Row 6
$ws.Range('D6').Value = "'0.6280"
$ws.Range('E6').Value = '  -1.55%  '

# Row 7
$ws.Range('D7').Value = "'1.004"
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('D8').Value = "'0.07534"
$ws.Range('E8').Value = '  -0.54%  '

# Row 9
$ws.Range('E9').Value = '  -0.12%  '

# Row 10
$ws.Range('D10').Value = "'23.37"
$ws.Range('E10').Value = '  +1.84%  '

# Row 11
$ws.Range('D11').Value = "'0.07729"
$ws.Range('E11').Value = '  -0.37%  '

# Row 12
$ws.Range('D12').Value = '1.869.10'
$ws.Range('E12').Value = '  +1.54%  '

# Row 13
$ws.Range('D13').Value = "'5.031"
$ws.Range('E13').Value = '  +0.24%  '

# Row 14
$ws.Range('D14').Value = "'0.6783"
$ws.Range('E14').Value = '  +0.69%  '

# Row 15
$ws.Range('D15').Value = "'83.30"
$ws.Range('E15').Value = '  -0.16%  '

# Row 16
$ws.Range('D16').Value = "'0.000009278"
$ws.Range('E16').Value = '  -3.17%  '

# Row 17
$ws.Range('D17').Value = "'5.991"
$ws.Range('E17').Value = '  -2.32%  '

# Row 18
$ws.Range('D18').Value = '29.242.37'
$ws.Range('E18').Value = '  +0.37%  '

# Row 19
$ws.Range('D19').Value = '2.096.02'
$ws.Range('E19').Value = '  +0.17%  '

# Row 20
$ws.Range('D20').Value = "'232.73"
$ws.Range('E20').Value = '  +2.41%  '

# Row 21
$ws.Range('D21').Value = "'12.75"
$ws.Range('E21').Value = '  +0.89%  '

# Row 22
$ws.Range('E22').Value = '  +0.28%  '

# Row 23
$ws.Range('D23').Value = "'7.196"
$ws.Range('E23').Value = '  -0.57%  '

# Row 24
$ws.Range('D24').Value = "'1.003"
$ws.Range('E24').Value = '  +0.14%  '

# Row 25
$ws.Range('D25').Value = "'160.52"
$ws.Range('E25').Value = '  -0.26%  '

# Row 26
$ws.Range('D26').Value = "'0.1399"
$ws.Range('E26').Value = '  -0.66%  '

# Row 27
$ws.Range('D27').Value = "'8.562"
$ws.Range('E27').Value = '  -0.09%  '

# Row 28
$ws.Range('D28').Value = "'17.98"
$ws.Range('E28').Value = '  -0.22%  '

# Row 29
$ws.Range('E29').Value = '  -0.32%  '

# Row 30
$ws.Range('D30').Value = "'4.198"
$ws.Range('E30').Value = '  +1.52%  '

# Row 31
$ws.Range('D31').Value = "'4.168"
$ws.Range('E31').Value = '  +1.87%  '

# Row 32
$ws.Range('D32').Value = "'0.05576"
$ws.Range('E32').Value = '  +3.08%  '

# Row 33
$ws.Range('E33').Value = '  +0.14%  '

# Row 34
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = "'0.7528"
$ws.Range('E34').Value = '  +0.31%  '

# Row 35
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = "'1.860"
$ws.Range('E35').Value = '  -0.44%  '

# Row 36
$ws.Range('E36').Value = '  +0.53%  '

# Row 37
$ws.Range('E37').Value = '  +0.07%  '

# Row 38
$ws.Range('D38').Value = '1.239.46'
$ws.Range('E38').Value = '  -0.23%  '

# Row 39
$ws.Range('D39').Value = "'2.775"
$ws.Range('E39').Value = '  +0.34%  '

# Row 40
$ws.Range('D40').Value = "'0.01795"
$ws.Range('E40').Value = '  +0.10%  '

# Row 41
$ws.Range('D41').Value = "'6.634"
$ws.Range('E41').Value = '  +0.21%  '

# Row 42
$ws.Range('D42').Value = "'0.9027"
$ws.Range('E42').Value = '  -0.46%  '

# Row 43
$ws.Range('E43').Value = '  +0.07%  '

# Row 44
$ws.Range('D44').Value = "'102.53"
$ws.Range('E44').Value = '  +0.15%  '

# Row 45
$ws.Range('D45').Value = '1.996.11'
$ws.Range('E45').Value = '  +0.25%  '

# Row 46
$ws.Range('D46').Value = "'66.75"
$ws.Range('E46').Value = '  +2.26%  '

# Row 47
$ws.Range('D47').Value = "'0.5109"
$ws.Range('E47').Value = '  -0.14%  '

# Row 48
$ws.Range('E48').Value = '  -3.22%  '

# Row 49
$ws.Range('D49').Value = "'0.4105"
$ws.Range('E49').Value = '  -0.17%  '

# Row 50
$ws.Range('D50').Value = "'9.143"
$ws.Range('E50').Value = '  +0.51%  '

# Row 51
$ws.Range('D51').Value = "'0.05849"
$ws.Range('E51').Value = '  +1.16%  '

# Reset style to remove quotePrefix formatting introduced by forcing text values
$resetCells = @("D5","D6","D7","D8","D10","D11","D13","D14","D15","D16","D17","D20","D21","D23","D24","D25","D26","D27","D28","D30","D31","D32","D34","D35","D39","D40","D41","D42","D44","D46","D47","D49","D50","D51")
foreach ($addr in $resetCells) {
    $ws.Range($addr).Style = "Normal"
}
